$d = $word.ActiveDocument

# 1) Merge "}% (from " + "figure" + ")" into a single run's text (cosmetic run-merge,
#    visible text unchanged): "}% (from figure)"
$d.Content.Find.Execute("}% (from figure)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "}% (from figure)", 2) | Out-Null

# 2) Merge "${" + "FPC" + "}% " into a single run: "${FPC}% "
$d.Content.Find.Execute("${FPC}% ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "${FPC}% ", 2) | Out-Null

# 3) Replace the hard-coded "$5,000" air tank price with the "${ATP}" placeholder.
$d.Content.Find.Execute("a new air tank will be `$5,000, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "a new air tank will be `${ATP}, ", 2) | Out-Null
